$d = $word.ActiveDocument

# Locate the "Night" section heading paragraph, then the first bullet
# paragraph after it whose text is exactly "250, 115, 30" (the Night
# section's second color entry, which sits right next to the _GoBack
# bookmark). We search from the heading forward so the identical-looking
# "250, 115, 30" bullet under "Daytime" is left untouched.
$nightIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Night*") {
        $nightIdx = $i
        break
    }
}

$targetIdx = -1
for ($i = $nightIdx; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "250, 115, 30") {
        $targetIdx = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIdx)
$pRange = $p.Range
$start = $pRange.Start

# The paragraph's visible text "250, 115, 30" is exactly 12 characters;
# replace just that span (leaving the trailing _GoBack bookmark, which
# sits right after it, untouched).
$old = $d.Range($start, $start + 12)
$old.Text = ""

# Insert the new color value as four separate runs, as produced by the
# edit: "250, " | "30" | ", " | "220"  ->  "250, 30, 220"
$ins = $d.Range($start, $start)
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = '<w:p ' + $ns + '>' +
       '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:t xml:space="preserve">250, </w:t></w:r>' +
       '<w:r><w:t>30</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
       '<w:r><w:t>220</w:t></w:r>' +
       '</w:p>'
$ins.InsertXML($xml)
